$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F with header "time_taken" (use header style from E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("F1").Style = $ws.Range("E1").Style

# Populate time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 13:42:23.270077"
$ws.Range("F3").Value = "2021-10-05 13:42:23.270090"
$ws.Range("F4").Value = "2021-10-05 13:42:23.270094"
$ws.Range("F5").Value = "2021-10-05 13:42:23.270097"
$ws.Range("F6").Value = "2021-10-05 13:42:23.270101"
$ws.Range("F7").Value = "2021-10-05 13:42:23.270104"
$ws.Range("F8").Value = "2021-10-05 13:42:23.270107"
$ws.Range("F9").Value = "2021-10-05 13:42:23.270110"
